$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7259.4375
$ws.Range("I19").Value = 20750.4
$ws.Range("J19").Value = 1127.1818
$ws.Range("K19").Value = 20750.4
$ws.Range("L19").Value = 1127.1818
$ws.Range("M19").Value = -20575.4
$ws.Range("N19").Value = -1477.1818
$ws.Range("H28").Value = 592.06665
$ws.Range("I28").Value = 578.4211
$ws.Range("J28").Value = 615.63635
$ws.Range("K28").Value = 578.4211
$ws.Range("L28").Value = 615.63635
$ws.Range("M28").Value = -93.42110000000002
$ws.Range("N28").Value = -1585.63635
$ws.Range("H33").Value = 6122.5293
$ws.Range("I33").Value = 12691.5
$ws.Range("J33").Value = 283.44446
$ws.Range("K33").Value = 12691.5
$ws.Range("L33").Value = 283.44446
$ws.Range("M33").Value = -12462.5
$ws.Range("N33").Value = -741.4444599999999
$ws.Range("H62").Value = 2834.6333
$ws.Range("I62").Value = 2475.1
$ws.Range("K62").Value = 2475.1
$ws.Range("M62").Value = -1851.1
$ws.Range("H65").Value = 2834.6333
$ws.Range("I65").Value = 2475.1
$ws.Range("K65").Value = 12375.5
$ws.Range("M65").Value = -9255.5
$ws.Range("H88").Value = 1189.7142
$ws.Range("J88").Value = 1236.909
$ws.Range("L88").Value = 1236.909
$ws.Range("N88").Value = -2048.909
$ws.Range("H91").Value = 1189.7142
$ws.Range("J91").Value = 1236.909
$ws.Range("L91").Value = 1236.909
$ws.Range("N91").Value = -4044.909
$ws.Range("H98").Value = 1002.0417
$ws.Range("I98").Value = 919.3889
$ws.Range("K98").Value = 919.3889
$ws.Range("M98").Value = 578.6111
$ws.Range("H122").Value = 1002.0417
$ws.Range("I122").Value = 919.3889
$ws.Range("K122").Value = 2758.1667
$ws.Range("M122").Value = -308.1667000000002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 50000
$ws.Range("J17").Value = 50000
$ws.Range("L17").Value = 50000
$ws.Range("N17").Value = -50346
$ws.Range("H45").Value = 1078.1
$ws.Range("I45").Value = 787.1053000000001
$ws.Range("K45").Value = 787.1053000000001
$ws.Range("M45").Value = -410.1053000000001
$ws.Range("H61").Value = 3764.5
$ws.Range("I61").Value = 1474.0834
$ws.Range("K61").Value = 1474.0834
$ws.Range("M61").Value = -1262.0834
$ws.Range("H63").Value = 18063.125
$ws.Range("I63").Value = 18063.125
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 18063.125
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -17377.125
$ws.Range("H66").Value = 18063.125
$ws.Range("I66").Value = 18063.125
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 90315.625
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -86883.625
$ws.Range("H97").Value = 1137.7567
$ws.Range("I97").Value = 1008.9375
$ws.Range("J97").Value = 1962.2
$ws.Range("K97").Value = 1008.9375
$ws.Range("L97").Value = 1962.2
$ws.Range("M97").Value = -512.9375
$ws.Range("N97").Value = -2954.2
$ws.Range("H110").Value = 60981.1
$ws.Range("I110").Value = 67533.11
$ws.Range("J110").Value = 2013
$ws.Range("K110").Value = 67533.11
$ws.Range("L110").Value = 2013
$ws.Range("M110").Value = -65488.11
$ws.Range("N110").Value = -6103
$ws.Range("H136").Value = 3764.5
$ws.Range("I136").Value = 1474.0834
$ws.Range("K136").Value = 4422.2502
$ws.Range("M136").Value = -1872.2502

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 101
$ws.Range("I8").Value = 101
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 101
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = 39
$ws.Range("H99").Value = 1259
$ws.Range("I99").Value = 880
$ws.Range("J99").Value = 1903.3
$ws.Range("K99").Value = 880
$ws.Range("L99").Value = 1903.3
$ws.Range("M99").Value = 618
$ws.Range("N99").Value = -4899.3
$ws.Range("H132").Value = 37250
$ws.Range("J132").Value = 37250
$ws.Range("L132").Value = 37250
$ws.Range("N132").Value = -47370

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1128.5333
$ws.Range("I16").Value = 1116.5
$ws.Range("K16").Value = 1116.5
$ws.Range("M16").Value = -829.5
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 1000
$ws.Range("M17").Value = -826
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = -900
$ws.Range("H86").Value = 2698.0908
$ws.Range("I86").Value = 3185.5
$ws.Range("J86").Value = 2113.2
$ws.Range("K86").Value = 3185.5
$ws.Range("L86").Value = 2113.2
$ws.Range("M86").Value = -2062.5
$ws.Range("N86").Value = -4359.2
$ws.Range("H89").Value = 2698.0908
$ws.Range("I89").Value = 3185.5
$ws.Range("J89").Value = 2113.2
$ws.Range("K89").Value = 15927.5
$ws.Range("L89").Value = 10566
$ws.Range("M89").Value = -10311.5
$ws.Range("N89").Value = -21798
$ws.Range("H99").Value = 100000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 100000
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").Value = 100000
$ws.Range("N99").Value = -102996
$ws.Range("H105").Value = 5701.5
$ws.Range("I105").Value = 5701.5
$ws.Range("K105").Value = 5701.5
$ws.Range("M105").Value = -3954.5
$ws.Range("H113").Value = 1128.5333
$ws.Range("I113").Value = 1116.5
$ws.Range("K113").Value = 1116.5
$ws.Range("M113").Value = 1053.5
$ws.Range("H126").Value = 100000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 100000
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").Value = 300000
$ws.Range("N126").Value = -304940
$ws.Range("H132").Value = 3355.8572
$ws.Range("I132").Value = 2913
$ws.Range("K132").Value = 8739
$ws.Range("M132").Value = -6209

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 1838
$ws.Range("J19").Value = 2225
$ws.Range("L19").Value = 6675
$ws.Range("N19").Value = -7023
$ws.Range("H68").Value = 3664281
$ws.Range("I68").Value = 8547423
$ws.Range("J68").Value = 1924.75
$ws.Range("K68").Value = 25642269
$ws.Range("L68").Value = 5774.25
$ws.Range("M68").Value = -25641458
$ws.Range("N68").Value = -7396.25
$ws.Range("H69").Value = 2275
$ws.Range("I69").Value = 925
$ws.Range("J69").Value = 2950
$ws.Range("K69").Value = 2775
$ws.Range("L69").Value = 8850
$ws.Range("M69").Value = -1964
$ws.Range("N69").Value = -10472
$ws.Range("H71").Value = 3664281
$ws.Range("I71").Value = 8547423
$ws.Range("J71").Value = 1924.75
$ws.Range("K71").Value = 76926807
$ws.Range("L71").Value = 17322.75
$ws.Range("M71").Value = -76922751
$ws.Range("N71").Value = -25434.75
$ws.Range("H72").Value = 2275
$ws.Range("I72").Value = 925
$ws.Range("J72").Value = 2950
$ws.Range("K72").Value = 8325
$ws.Range("L72").Value = 26550
$ws.Range("M72").Value = -4269
$ws.Range("N72").Value = -34662
$ws.Range("H92").Value = 575.5
$ws.Range("I92").Value = 434
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 1302
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -54
$ws.Range("N92").Value = -5496
$ws.Range("H109").Value = 2345
$ws.Range("I109").Value = 312.77777
$ws.Range("K109").Value = 938.33331
$ws.Range("M109").Value = 101.66669
$ws.Range("H115").Value = 4499.273
$ws.Range("J115").Value = 5199.1113
$ws.Range("L115").Value = 15597.3339
$ws.Range("N115").Value = -17947.3339
$ws.Range("H119").Value = 5128.923
$ws.Range("I119").Value = 3542.5715
$ws.Range("J119").Value = 6979.6665
$ws.Range("K119").Value = 10627.7145
$ws.Range("L119").Value = 20938.9995
$ws.Range("M119").Value = -5789.7145
$ws.Range("N119").Value = -30614.9995
$ws.Range("H131").Value = 2204.3777
$ws.Range("J131").Value = 2241.8977
$ws.Range("L131").Value = 6725.6931
$ws.Range("N131").Value = -16805.6931

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 50000
$ws.Range("J6").Value = 50000
$ws.Range("L6").Value = 50000
$ws.Range("N6").Value = -50226
$ws.Range("H16").Value = 50000
$ws.Range("J16").Value = 50000
$ws.Range("L16").Value = 50000
$ws.Range("N16").Value = -50500
$ws.Range("H80").Value = 3926.375
$ws.Range("I80").Value = 4667.3335
$ws.Range("J80").Value = 1703.5
$ws.Range("K80").Value = 4667.3335
$ws.Range("L80").Value = 1703.5
$ws.Range("M80").Value = -3669.3335
$ws.Range("N80").Value = -3699.5
$ws.Range("H83").Value = 3926.375
$ws.Range("I83").Value = 4667.3335
$ws.Range("J83").Value = 1703.5
$ws.Range("K83").Value = 23336.6675
$ws.Range("L83").Value = 8517.5
$ws.Range("M83").Value = -18344.6675
$ws.Range("N83").Value = -18501.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 2402
$ws.Range("I19").Value = 300
$ws.Range("J19").Value = 2927.5
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 2927.5
$ws.Range("M19").Value = -130
$ws.Range("N19").Value = -3267.5
$ws.Range("H61").Value = 1910.2142
$ws.Range("I61").Value = 1113.5
$ws.Range("J61").Value = 2972.5
$ws.Range("K61").Value = 1113.5
$ws.Range("L61").Value = 2972.5
$ws.Range("M61").Value = -911.5
$ws.Range("N61").Value = -3376.5
$ws.Range("H113").Value = 1910.2142
$ws.Range("I113").Value = 1113.5
$ws.Range("J113").Value = 2972.5
$ws.Range("K113").Value = 1113.5
$ws.Range("L113").Value = 2972.5
$ws.Range("M113").Value = 1056.5
$ws.Range("N113").Value = -7312.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1420.6875
$ws.Range("I96").Value = 1221.9
$ws.Range("J96").Value = 1752
$ws.Range("K96").Value = 1221.9
$ws.Range("L96").Value = 1752
$ws.Range("M96").Value = 151.0999999999999
$ws.Range("N96").Value = -4498
$ws.Range("H132").Value = 2458.8147
$ws.Range("I132").Value = 2021.3043
$ws.Range("K132").Value = 6063.9129
$ws.Range("M132").Value = -3533.9129
